# Generate Report for Handoff
# Replace the old handoff-package identifiers / timestamps with the new ones
# produced by this handoff run, on all three worksheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$oldGuid = "4db06b64-63f1-4790-a880-b2d4e60f3865"
$newGuid = "bf20ca34-e3a8-4637-adfa-dad0620a5e1d"

$oldHash = "0ae0b0276c1414e9a737df68b7538097aa47b33b"
$newHash = "24482dccac05babb4228417968adf195803df08e"

$newMdName  = "$newGuid.md"
$newZhName  = "$newGuid.$newHash.zh-cn.xlf"
$newDeName  = "$newGuid.$newHash.de-de.xlf"

$newHandoffDate = "2016-03-22 00:55:44"
$newZhDatetime  = "2016-03-22 00:55:40"

# Original (pre-edit) hyperlink target URLs, taken from the workbook as
# authored. These underlying link targets are untouched by this commit
# (only the cell text / hyperlink display text is refreshed for the new
# handoff run), so we keep them pointing at the original package paths.
# (Hyperlink.Address/SubAddress/ScreenTip read back as blank through this
# engine's object model, and deleting a hyperlink anywhere on a worksheet
# clears the whole sheet's Hyperlinks collection, so every link on a sheet
# has to be re-added together with its unchanged address and its correct,
# final display text.)
$oldMdName = "$oldGuid.md"
$oldZhName = "$oldGuid.$oldHash.zh-cn.xlf"
$oldDeName = "$oldGuid.$oldHash.de-de.xlf"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/e338b2bf903720b0a1bcdc585e1455a57a1438c6/e2e/$oldMdName"
$zhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/db278d7ff44ada3dc28467fc1240bb36e9e5ec0e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZhName"
$deUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/025ccaa10a717c53d044b3de2833634fa9fa04e4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDeName"

function Reset-SheetHyperlinks($ws, $links) {
    # $links: array of hashtables @{ Cell = "A2"; Address = "..."; Display = "..." }
    if ($ws.Hyperlinks.Count -gt 0) {
        $ws.Range($links[0].Cell).Hyperlinks.Delete()
    }
    foreach ($l in $links) {
        $ws.Hyperlinks.Add($ws.Range($l.Cell), $l.Address, "", "", $l.Display) | Out-Null
    }
}

# ---------- Overview sheet ----------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $newHandoffDate
Reset-SheetHyperlinks $wsOverview @(
    @{ Cell = "A2"; Address = $mdUrl; Display = $newMdName }
)

# ---------- zh-cn sheet ----------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhName
$wsZh.Range("E2").Value = $newZhDatetime
Reset-SheetHyperlinks $wsZh @(
    @{ Cell = "A2"; Address = $mdUrl; Display = $newMdName },
    @{ Cell = "D2"; Address = $zhUrl; Display = $newZhName }
)

# ---------- de-de sheet ----------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeName
Reset-SheetHyperlinks $wsDe @(
    @{ Cell = "A2"; Address = $mdUrl; Display = $newMdName },
    @{ Cell = "D2"; Address = $deUrl; Display = $newDeName }
)

Write-Output "Report regenerated for handoff."
